$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new year columns (K = 2021, L = 2022) are being appended to the table
# that currently runs through column J. For every existing column-J cell we
# copy its full formatting (number format / font / borders / alignment) onto
# the new K/L cells via PasteSpecial, then fix up alignment / values as
# needed, mirroring how the existing sheet was built column-by-column.

function Copy-RowFormat([string]$src, [string]$dstRange) {
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dstRange).PasteSpecial(-4122) | Out-Null
}

# Row 3 - bottom border separator row (same style as existing D3:J3, s="3")
Copy-RowFormat "J3" "K3:L3"

# Row 4 - bold year header row (same style as existing D4:J4, s="10")
Copy-RowFormat "J4" "K4:L4"
$ws.Range("K4").Value = 2021
$ws.Range("L4").Value = 2022

# Row 5 - first data row, has a top border (like J5, s="14") but the new
# columns drop the "horizontal=right" alignment (becomes new style s="27")
Copy-RowFormat "J5" "K5:L5"
$ws.Range("K5:L5").HorizontalAlignment = 1
$ws.Range("K5").Value = 272.6
$ws.Range("L5").Value = 292.19961890663211

# Row 6 - empty subtotal-label row (like J6, s="17" -> new s="28")
Copy-RowFormat "J6" "K6:L6"
$ws.Range("K6:L6").HorizontalAlignment = 1

# Row 7
Copy-RowFormat "J7" "K7:L7"
$ws.Range("K7:L7").HorizontalAlignment = 1
$ws.Range("K7").Value = 98.1
$ws.Range("L7").Value = 99.522498012012946

# Row 8
Copy-RowFormat "J8" "K8:L8"
$ws.Range("K8:L8").HorizontalAlignment = 1
$ws.Range("K8").Value = 174.5
$ws.Range("L8").Value = 192.67712089461918

# Row 9 - empty subtotal-label row
Copy-RowFormat "J9" "K9:L9"
$ws.Range("K9:L9").HorizontalAlignment = 1

# Row 10
Copy-RowFormat "J10" "K10:L10"
$ws.Range("K10:L10").HorizontalAlignment = 1
$ws.Range("K10").Value = 75.6
$ws.Range("L10").Value = 88.011952928467494

# Row 11
Copy-RowFormat "J11" "K11:L11"
$ws.Range("K11:L11").HorizontalAlignment = 1
$ws.Range("K11").Value = 55.5
$ws.Range("L11").Value = 56.919430260413804

# Row 12 - bottom row with bottom border (like J12, s="19" -> new s="29")
Copy-RowFormat "J12" "K12:L12"
$ws.Range("K12:L12").HorizontalAlignment = 1
$ws.Range("K12").Value = 24.9
$ws.Range("L12").Value = 24.176373211436804

$ws.Range("N5").Select() | Out-Null
